# Apply updated "K" column (G) values to Sheet1, rows 2-13.
# These values reflect a regeneration of save_data (K replacing Strike#,
# recalculated std/mean, s_vals) — only column G changes per the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newK = @{
    2  = 5
    3  = 4
    4  = 3
    5  = 4
    6  = 3
    7  = 5
    8  = 2
    9  = 3
    10 = 1
    11 = 7
    12 = 3
    13 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
